$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 16 with the new status entry (20/1/2021) ---

# A16: date text, same style as A12/A13/A14/A15 (left aligned, no wrap)
$ws.Range("A16").Value = "20/1/2021"
$ws.Range("A16").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A16").WrapText = $false

# D16: related files - written before B16 so the shared-string table order
# matches the target (D16's string becomes index 26, B16's becomes index 27)
$ws.Range("D16").Value = "Flipping_Bits.txt`nBitwise_Operators.txt"
$ws.Range("D16").WrapText = $true

# B16: long status notes, wrap text like B15
$ws.Range("B16").Value = "1. Added few more GIT commands in the cheat sheet`n2. Attended syncup meeting with Srivalli`n3. Completed bitwise operators in C and did 2 hacker rank programs on bitwise operators`n4. Completed defining and accessing members of structures`n5. Understand how memory is allocated for structures`n6. Completed how to define and access members of unions and understood memory allocation of members`n7. Understood bitfields `n8. Attended meeting ""Softwate testing session by Srinivasa"""
$ws.Range("B16").WrapText = $true

# Row 16 needs to grow tall enough to show the whole note (target height 225)
$ws.Rows(16).RowHeight = 225

# Column D is no longer best-fit; it now has a fixed, wider width
$ws.Range("D1").EntireColumn.ColumnWidth = 19.8

# --- Update the view: scroll down a row and move the active selection ---
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
